# Update for Vehicle Renewal & Login For Corporate
#
# - AlternativeDataRegisterVehicle: view was scrolled (topLeftCell="B1");
#   reset the zoom/view so it no longer shows a custom top-left cell.
# - vehicleRenewalTestData: selection moved to N13, and a few test-data
#   rows were updated (VehicleWeightFrom/To values + toRun flags flipped
#   to True for the first two scenarios).

$wb = $excel.ActiveWorkbook
$origActiveSheetName = $excel.ActiveSheet.Name

# --- AlternativeDataRegisterVehicle: clear scrolled view / reset zoom ---
$wsAlt = $wb.Worksheets.Item("AlternativeDataRegisterVehicle")
$wsAlt.Activate() | Out-Null
$altWindow = $excel.ActiveWindow
$altWindow.Zoom = 100

# --- vehicleRenewalTestData: selection + data updates ---
$wsRenewal = $wb.Worksheets.Item("vehicleRenewalTestData")
$wsRenewal.Activate() | Out-Null
$wsRenewal.Range("N13").Select() | Out-Null

# Row 2 (TC 001): VehicleWeightFrom 3000 -> 120, VehicleWeightTo 4000 -> 300,
# toRun False -> True. The leading apostrophe keeps these as text entries
# (column is formatted as Text) instead of being coerced to numbers.
$wsRenewal.Range("B2").Value = "'120"
$wsRenewal.Range("C2").Value = "300"
$wsRenewal.Range("G4").Copy($wsRenewal.Range("G2"))

# Row 3 (TC 002): VehicleWeightFrom 3000 -> 3100, toRun False -> True.
$wsRenewal.Range("B3").Value = "'3100"
$wsRenewal.Range("G4").Copy($wsRenewal.Range("G3"))

# Restore original active sheet
$wb.Worksheets.Item($origActiveSheetName).Activate() | Out-Null
